$d = $word.ActiveDocument

# Locate the paragraph containing the tagline that needs updating.
$oldText = "* I make things - things that work. *"
$findRng = $d.Content
$found = $findRng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text to replace"
}

$para = $findRng.Paragraphs(1)
$paraRange = $para.Range

$lsq = [char]0x2018
$rsq = [char]0x2019

# New italic runs that replace the old, non-italic "* I make things ... *" run.
# (The leading "Platform - Tools - Infrastructure - Security" run and the
# single space run that follow it are left exactly as they were.)
$newRuns = @(
    "I make things - things that work; and by",
    " ",
    "$lsq",
    "work",
    "$rsq",
    " ",
    "I mean work superlatively."
)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function XmlEscape([string]$s) {
    $s = $s.Replace("&", "&amp;")
    $s = $s.Replace("<", "&lt;")
    $s = $s.Replace(">", "&gt;")
    return $s
}

$runsXml = ""
$runsXml += "<w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space='preserve'>" + (XmlEscape "Platform - Tools - Infrastructure - Security") + "</w:t></w:r>"
$runsXml += "<w:r><w:t xml:space='preserve'> </w:t></w:r>"
foreach ($t in $newRuns) {
    $runsXml += "<w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space='preserve'>" + (XmlEscape $t) + "</w:t></w:r>"
}

$paraXml = "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='FirstParagraph'/></w:pPr>$runsXml</w:p>"

$paraRange.InsertXML($paraXml)
